# Insert a new data row at Excel row 268 (shifts existing rows 268:391 down
# to 269:392, which matches the rest of the rows "moving down by one" seen
# in the diff). The new row inherits formatting (incl. the date number
# format in column D) from the row above, same as Excel's default Insert
# behavior.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows(268).Insert()

# Populate the newly inserted row 268 with its data.
$ws.Cells.Item(268, 1).Value = 10
$ws.Cells.Item(268, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(268, 3).Value = "La Araucanía"
$ws.Cells.Item(268, 4).Value = 44875
$ws.Cells.Item(268, 5).Value = 9
$ws.Cells.Item(268, 6).Value = 100112017
$ws.Cells.Item(268, 7).Value = "Apio"
$ws.Cells.Item(268, 8).Value = "Americana (o)"
$ws.Cells.Item(268, 9).Value = "Primera"
$ws.Cells.Item(268, 10).Value = 110
$ws.Cells.Item(268, 11).Value = 10000
$ws.Cells.Item(268, 12).Value = 10000
$ws.Cells.Item(268, 13).Value = 10000
$ws.Cells.Item(268, 14).Value = "`$/docena de matas"
$ws.Cells.Item(268, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(268, 16).Value = 1667
$ws.Cells.Item(268, 17).Value = 6
$ws.Cells.Item(268, 18).Value = "Hortaliza"
